$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (single-dot decimal-looking price strings) get forced to Text format first,
# then restored to the default (General) style so only the value text changes.

$ws.Range("D2").Value = '27.069.95'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.680.07'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  -2.94%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +1.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.39'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.65%  '
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("D12").Value = '1.916.73'
$ws.Range("D13").Value = '1.661.02'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.535'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").Value = '27.066.15'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("E18").Value = '  +2.14%  '
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("E20").Value = '  +0.58%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.34%  '
$ws.Range("E23").Value = '  +1.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.97%  '
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = '1.546.18'
$ws.Range("E33").Value = '  +6.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("E35").Value = '  +4.62%  '
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.589'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("E40").Value = '  +7.15%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").Value = '1.821.72'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.781'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.41%  '
